$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status went from "Ready for handoff" to "In Translation" for both
# e2e files, in both the per-language sheets and the rolled-up Overview
# sheet (columns E/F there mirror the zh-cn/de-de Status column).
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

$ovw.Range("E2").Value = "In Translation"
$ovw.Range("E3").Value = "In Translation"
$ovw.Range("F2").Value = "In Translation"
$ovw.Range("F3").Value = "In Translation"

# The Status text got shorter, so the Status columns (and their mirrors
# on the Overview sheet) shrink when re-fit to the new content.
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
$ovw.Range("E1:F1").EntireColumn.AutoFit()
